# Time Log.xlsx - add the missing time-entry row (row 113) for 2014-11-13,
# an 11-minute "Coding" session with a 5-minute interruption, then move the
# active selection to B114 (ready for the next entry), matching the
# "ListView scrolls to currently used background." commit.

# NOTE on ordering: the Interruption (column D) cell is deliberately written
# before the Start/Stop time cells (B, C). The shared formula in column E
# reads B, C and D together; writing D first ensures the engine's
# dependency tracking picks up the interruption minutes when it evaluates
# the E113 shared formula (writing D last leaves E113 computed as if D
# were still 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Date (serial 41956 == 2014-11-13)
$ws.Cells.Item(113, 1).Value = 41956
# Interruption (minutes) - set before the time columns, see note above.
$ws.Cells.Item(113, 4).Value = 5
# Start time (15:27) / Stop time (15:58), stored as fraction-of-day
$ws.Cells.Item(113, 2).Value = 0.64374999999999993
$ws.Cells.Item(113, 3).Value = 0.66527777777777775
# Activity / category
$ws.Cells.Item(113, 6).Value = "Coding"

# Recalculate so E113 (shared formula), the Sheet1 total, and the Sheet2
# SUMIF/percentage table all pick up the new row.
$excel.Calculate()

# Move the selection on to the next empty row, column B, as in the diff.
$ws.Range("B114").Select()
